# Fill in the computed Retention ratio and Answer Recall Average (ARA)
# values that were left blank in the results table.
#
# The document contains a single table whose last column holds numeric
# results for various metrics. The following rows currently have an
# empty value cell that needs to be populated with a bold run:
#   Row 24 ("Ratio")                       -> 0.2857
#   Row 44 ("Answer Recall Lenient (ARL)") -> 0.3571
#   Row 45 ("Answer Recall Strict (ARS)")  -> 0.1428
#   Row 46 ("Answer Recall Average (ARA)") -> 0.2499

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-ResultValue($rowIndex, $value) {
    $row = $t.Rows.Item($rowIndex)
    $cell = $row.Cells.Item($row.Cells.Count)
    $cell.Range.Text = $value
    $valueRange = $cell.Range
    $valueRange.Font.Bold = 1
    $valueRange.Font.Size = 12
    $valueRange.Font.SizeBi = 12
}

Set-ResultValue 24 "0.2857"
Set-ResultValue 44 "0.3571"
Set-ResultValue 45 "0.1428"
Set-ResultValue 46 "0.2499"
